$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source sheet stores every cell (coin name, link, price, % change) as plain text.
# Most replacement strings (URLs, "%"-suffixed deltas, multi-dot prices like "59.313.63")
# can never be parsed as a number, so a plain .Value assignment round-trips them as text
# exactly like the original inlineStr cells. A few new price strings DO look like plain
# numbers (e.g. "6.45", "0.0000134") -- for just those we flip the cell to Text format
# first so Excel keeps the literal digits/trailing zeros instead of normalising them into
# a float.

$ws.Range("D2").Value = "59.313.63"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").Value = "2.605.02"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.00"
$ws.Range("E5").Value = "  +4.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.69"
$ws.Range("E6").Value = "  +1.68%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.45"
$ws.Range("E9").Value = "  -1.33%  "
$ws.Range("E10").Value = "  +2.55%  "
$ws.Range("E11").Value = "  +1.54%  "
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").Value = "3.060.70"
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("D14").Value = "59.247.71"
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.62"
$ws.Range("E15").Value = "  +1.12%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.622.63"
$ws.Range("E16").Value = "  +1.26%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000134"
$ws.Range("E17").Value = "  +1.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "341.24"
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("E19").Value = "  +1.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.15"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("E21").Value = "  -1.24%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.71"
$ws.Range("E23").Value = "  +2.46%  "
$ws.Range("E24").Value = "  +1.45%  "
$ws.Range("E25").Value = "  -1.54%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  +3.06%  "
$ws.Range("D28").Value = "0.0₃0751"
$ws.Range("E28").Value = "  +4.78%  "
$ws.Range("E30").Value = "  +8.08%  "
$ws.Range("E31").Value = "  -2.25%  "
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.62"
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.98"
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.13"
$ws.Range("E35").Value = "  +0.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "37.17"
$ws.Range("E36").Value = "  +2.32%  "
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.836"
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.829"
$ws.Range("E39").Value = "  +1.96%  "
$ws.Range("E40").Value = "  +2.09%  "
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "275.23"
$ws.Range("E42").Value = "  +0.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.599"
$ws.Range("E43").Value = "  +1.98%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0955"
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0523"
$ws.Range("E46").Value = "  +0.83%  "
$ws.Range("D47").Value = "1.953.91"
$ws.Range("E47").Value = "  -1.30%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.58"
$ws.Range("E48").Value = "  +4.13%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0224"
$ws.Range("E49").Value = "  +1.65%  "
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("E51").Value = "  -0.63%  "
